# Ajout draft mapping f595a2bd5e53be80aa00972cfd76eee4a5f7087b
#
#  - bump the IG build Date on the "Metadata" sheet
#  - add a new mapping column to the "Elements" sheet:
#      "Mapping: Spécification métier vers l'extension ROR ConfidentialityLevel"
#    with the draft mapping target "niveauConfidentialite" recorded against
#    the Extension.value[x] row (the only row with a mapping so far).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet: refresh the generation Date (row 8, column B)
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# ---------------------------------------------------------------------------
# Elements sheet: append the new "Mapping: ..." column (column 38 / AL),
# right after the existing "Mapping: RIM Mapping" column (AK).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Header cell: copy the formatting of the neighbouring mapping header so the
# new column keeps the same bold/banner style, then set its text.
$ws.Range("AK1").Copy($ws.Range("AL1"))
$ws.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR ConfidentialityLevel"

# The new column is blank for every element row except the
# Extension.value[x] row (row 6). Copy an existing blank data cell's
# formatting into the new column for those rows so they match the sheet's
# normal (bordered/wrapped) cell style.
$ws.Range("D2").Copy($ws.Range("AL2"))
$ws.Range("D2").Copy($ws.Range("AL3"))
$ws.Range("D2").Copy($ws.Range("AL4"))
$ws.Range("D2").Copy($ws.Range("AL5"))

# Extension.value[x] (row 6) now maps to the business term
# "niveauConfidentialite" - copy the row's existing data-cell style first,
# then write the new mapping value.
$ws.Range("AK6").Copy($ws.Range("AL6"))
$ws.Range("AL6").Value = "niveauConfidentialite"

# Match the column's best-fit width to the new (wide) header text.
$ws.Columns.Item(38).ColumnWidth = 75
